# Remove the trailing "site footer" block that was appended to the end of
# the bibliography: a blank separator paragraph, the "Ver no Jupiter ..."
# line and the "(c) 2020 ... Creative Commons Attribution" line.
#
# The paragraph right after the removed block (another blank "Normal"
# paragraph) and everything before the separator stay untouched.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startPara = $p
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    # Also drop the blank separator paragraph immediately preceding the
    # "Ver no Jupiter" line so no empty line is left behind.
    $deleteFrom = $startPara
    $prevPara = $startPara.Previous()
    if ($prevPara -ne $null -and $prevPara.Range.Text -eq "`r") {
        $deleteFrom = $prevPara
    }

    $r = $d.Range($deleteFrom.Range.Start, $endPara.Range.End)
    $r.Delete()
}
